$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.602.30'
$ws.Range("E2").Value = '  +1.17%  '

$ws.Range("D3").Value = '3.395.95'
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.19'
$ws.Range("E5").Value = '  +1.02%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.57'
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.72'
$ws.Range("E9").Value = '  +2.46%  '

$ws.Range("E10").Value = '  -0.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.389'
$ws.Range("E11").Value = '  -1.55%  '

$ws.Range("D12").Value = '3.975.65'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("E13").Value = '  +0.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.43'
$ws.Range("E14").Value = '  +0.72%  '

$ws.Range("D15").Value = '3.386.90'
$ws.Range("E15").Value = '  -0.30%  '

$ws.Range("E16").Value = '  +0.25%  '

$ws.Range("D17").Value = '61.586.10'
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.16'
$ws.Range("E18").Value = '  -0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.72'
$ws.Range("E19").Value = '  -1.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.00'
$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.80'
$ws.Range("E21").Value = '  +1.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.61'
$ws.Range("E22").Value = '  +1.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.555'
$ws.Range("E23").Value = '  -0.58%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("E25").Value = '  -3.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.194'
$ws.Range("E26").Value = '  +9.08%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  -1.15%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.05'
$ws.Range("E29").Value = '  +0.95%  '

$ws.Range("E30").Value = '  +0.47%  '

$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("E32").Value = '  -3.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.43'
$ws.Range("E33").Value = '  -0.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.96'
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.73'
$ws.Range("E35").Value = '  +0.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.06'
$ws.Range("E36").Value = '  +1.76%  '

$ws.Range("D37").Value = '3.430.47'
$ws.Range("E37").Value = '  +0.31%  '

$ws.Range("E38").Value = '  -0.54%  '

$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.97'
$ws.Range("E40").Value = '  -5.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.782'
$ws.Range("E41").Value = '  +0.18%  '

$ws.Range("E42").Value = '  +0.24%  '

$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.16'
$ws.Range("E44").Value = '  +1.50%  '

$ws.Range("D45").Value = '2.474.93'
$ws.Range("E45").Value = '  -0.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.10'
$ws.Range("E46").Value = '  +0.13%  '

$ws.Range("E47").Value = '  -1.90%  '

$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0264'
$ws.Range("E49").Value = '  -1.51%  '

$ws.Range("E50").Value = '  -0.79%  '

$ws.Range("E51").Value = '  -1.15%  '
